$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.591.32'
$ws.Range("E2").Value = '  -1.74%  '
$ws.Range("D3").Value = '2.901.50'
$ws.Range("E3").Value = '  -2.62%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '526.37'
$ws.Range("E5").Value = '  -2.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.94'
$ws.Range("E6").Value = '  -5.67%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.547'
$ws.Range("E8").Value = '  -3.59%  '
$ws.Range("D9").Value = '2.907.37'
$ws.Range("E9").Value = '  -2.82%  '
$ws.Range("E10").Value = '  -5.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.98'
$ws.Range("E11").Value = '  -2.71%  '
$ws.Range("E12").Value = '  -2.86%  '
$ws.Range("D13").Value = '3.405.12'
$ws.Range("E13").Value = '  -2.74%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.129'
$ws.Range("E14").Value = '  +3.40%  '
$ws.Range("D15").Value = '60.574.57'
$ws.Range("E15").Value = '  -1.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.48'
$ws.Range("E16").Value = '  -5.68%  '
$ws.Range("D17").Value = '2.909.02'
$ws.Range("E17").Value = '  -2.41%  '
$ws.Range("E18").Value = '  -4.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.96'
$ws.Range("E19").Value = '  -3.97%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.56'
$ws.Range("E20").Value = '  -4.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '350.12'
$ws.Range("E21").Value = '  -8.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.49'
$ws.Range("E22").Value = '  -3.09%  '
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.72'
$ws.Range("E24").Value = '  +1.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.72'
$ws.Range("E25").Value = '  -1.76%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.450'
$ws.Range("E26").Value = '  -4.53%  '
$ws.Range("E27").Value = '  -6.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("E29").Value = '  -4.52%  '
$ws.Range("D30").Value = '0.0₃0854'
$ws.Range("E30").Value = '  -9.00%  '
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("E32").Value = '  -2.88%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.48'
$ws.Range("E33").Value = '  -5.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '151.80'
$ws.Range("E34").Value = '  -4.92%  '
$ws.Range("E35").Value = '  -6.06%  '
$ws.Range("E36").Value = '  -6.40%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.991'
$ws.Range("E37").Value = '  -7.09%  '
$ws.Range("E38").Value = '  -6.34%  '
$ws.Range("E39").Value = '  +0.09%  '
$ws.Range("E40").Value = '  -5.33%  '
$ws.Range("E41").Value = '  -5.04%  '
$ws.Range("D42").Value = '2.288.52'
$ws.Range("E42").Value = '  -5.40%  '
$ws.Range("E43").Value = '  -3.48%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0577'
$ws.Range("E44").Value = '  -2.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.31'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.92'
$ws.Range("E47").Value = '  -3.87%  '
$ws.Range("E48").Value = '  -3.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.33'
$ws.Range("E49").Value = '  -0.97%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0915'
$ws.Range("E50").Value = '  -3.94%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.26'
$ws.Range("E51").Value = '  -7.52%  '
